# feat: unical ids of dialog
#
# Prefixes every "replic-id" (column C) and every "replic-position" (column D)
# with "test." (the sheet's own source/root id), and prefixes every id that is
# referenced inside the "replic-includes" list (column H) the same way.
# Rows 70-72 already carry ids that start with "test." (role definitions for
# the "test" dialog source) and get a second "test." prefix, matching the
# target data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters -> numbers used on this sheet.
$colC = 3   # replic-id
$colD = 4   # replic-position
$colH = 8   # replic-includes

$lastRow = 72

# --- D2: empty "replic-position" on the root row becomes "test." ---------
$ws.Cells.Item(2, $colD).Value2 = "test."

# --- Column C: prefix every replic-id except the root row (row 2, "test") -
for ($r = 3; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colC)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $cell.Value2 = "test." + $val
    }
}

# --- Column D: prefix every replic-position except row 3 ------------------
# (row 3's position points at the root row, whose id stays "test" unprefixed)
# Rows 70-72 are role definitions whose position ("test") is left untouched.
for ($r = 4; $r -le 69; $r++) {
    $cell = $ws.Cells.Item($r, $colD)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $cell.Value2 = "test." + $val
    }
}

# --- Column H: prefix every quoted id inside the includes list ------------
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colH)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $found = [regex]::Matches($val, "'([^']*)'")
        if ($found.Count -gt 0) {
            $parts = @()
            foreach ($m in $found) {
                $parts += "'test." + $m.Groups[1].Value + "'"
            }
            $cell.Value2 = "[" + ($parts -join ", ") + "]"
        }
    }
}
